$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold an exact literal string (avoids Excel auto-converting
# numeric-looking text like "215.30" into a number and losing formatting),
# while keeping the cells style identical to its original (no style).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.139.95"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.643.69"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +1.01%  "
Set-TextValue $ws.Range("D5") "215.30"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +1.05%  "
Set-TextValue $ws.Range("D8") "0.249"
$ws.Range("E8").Value = "  -3.00%  "
Set-TextValue $ws.Range("D9") "0.0621"
$ws.Range("E9").Value = "  -2.75%  "
Set-TextValue $ws.Range("D10") "18.60"
$ws.Range("E10").Value = "  -5.08%  "
Set-TextValue $ws.Range("D11") "0.0794"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.724.20"
$ws.Range("E13").Value = "  -1.83%  "
Set-TextValue $ws.Range("D14") "0.528"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "62.38"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.141.04"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "0.0₃0745"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("E18").Value = "  +0.99%  "
Set-TextValue $ws.Range("D19") "190.08"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  -2.43%  "
Set-TextValue $ws.Range("D21") "9.52"
$ws.Range("E21").Value = "  -4.14%  "
Set-TextValue $ws.Range("D22") "6.08"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("E23").Value = "  +1.54%  "
Set-TextValue $ws.Range("D24") "143.69"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +0.39%  "
Set-TextValue $ws.Range("D26") "1.77"
$ws.Range("E26").Value = "  -1.61%  "
Set-TextValue $ws.Range("D27") "6.71"
$ws.Range("E27").Value = "  -1.97%  "
Set-TextValue $ws.Range("D28") "15.14"
$ws.Range("E28").Value = "  -2.75%  "
Set-TextValue $ws.Range("D29") "1.24"
$ws.Range("E29").Value = "  -0.14%  "
Set-TextValue $ws.Range("D30") "0.0479"
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("E31").Value = "  -2.60%  "
Set-TextValue $ws.Range("D32") "3.14"
$ws.Range("E32").Value = "  -5.24%  "
Set-TextValue $ws.Range("D33") "2.45"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("E34").Value = "  -1.85%  "
Set-TextValue $ws.Range("D35") "0.876"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").Value = "1.125.06"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.518"
$ws.Range("E38").Value = "  -5.16%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.0154"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D40") "98.39"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "0.789"
$ws.Range("E41").Value = "  -1.07%  "
Set-TextValue $ws.Range("D42") "5.26"
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("D43").Value = "0.0₆0113"
$ws.Range("E43").Value = "  -1.54%  "
Set-TextValue $ws.Range("D44") "55.27"
$ws.Range("E44").Value = "  -2.64%  "
Set-TextValue $ws.Range("D45") "0.0521"
$ws.Range("E45").Value = "  -1.79%  "
Set-TextValue $ws.Range("D46") "1.47"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D47") "0.418"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "7.58"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("E49").Value = "  +0.87%  "
Set-TextValue $ws.Range("D50") "0.0922"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("E51").Value = "  -1.32%  "
